# Apply the "Added back transition col" edit:
#  - K5 utterance text changes from "Transition COL" to "I see."
#  - K5 keeps the same visual formatting (Arial 11, black, not bold,
#    yellow fill / thin border inherited from the row, general alignment)
#  - the active selection moves to K13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("K5")
$cell.Value = "I see."

# Re-assert the cell's formatting explicitly (size 11, regular weight,
# pure black Arial text) so the font/style used for this cell matches
# the one described for the edited row.
$cell.Font.Name = "Arial"
$cell.Font.Size = 11
$cell.Font.Bold = $false
$cell.Font.Color = 0
$cell.HorizontalAlignment = 1

# Move the active selection, matching the workbook's recorded cursor
# position after the edit.
$ws.Range("K13").Select() | Out-Null
